$d = $word.ActiveDocument

# Find the anchor paragraph: "Admin specific věci – to be done"
$anchorIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "Admin specific*to be done*") {
        $anchorIndex = $i
        break
    }
}

$lines = @(
    "Resetování hesel uživatelů",
    "",
    "",
    "",
    "",
    "Zobrazení od účetní",
    "",
    "Customer:",
    "First name",
    "Last name",
    "Title",
    "Address",
    "Postal code",
    "City",
    "Email",
    "Telephone",
    "",
    "Je loyal – bool RO",
    "Username – RO (pouze pokud je loyal, jinak schovat pole)",
    "Account id – RO",
    "Base address – RO",
    "RO – Read Only",
    ""
)

$idx = $anchorIndex
foreach ($line in $lines) {
    $p = $d.Paragraphs.Item($idx)
    $p.Range.InsertParagraphAfter()
    $idx = $idx + 1
    if ($line -ne "") {
        $newp = $d.Paragraphs.Item($idx)
        $newp.Range.Text = $line
    }
}
